$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "A"
$ws.Range("B8").Formula = "=(C2+C3)/2*C4"

$ws.Range("A10").Value = "Atotal"
$ws.Range("B10").Value = "Qtotal"

$ws.Range("A11").Formula = "=B8"
$ws.Range("B11").Formula = "=B7"

$ws.Range("A11:B11").Select()
